# Applies the "Modifing the lab_test table and add permanent value too it."
# edit: fills in a new backlog row (row 8) on the burndown sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New backlog item in row 8. Shared-string table order matters: the
# "Modifing..." text must be interned before "Redesign the database." so
# the new entries land at sharedStrings indices 25/26 respectively.
$ws.Range("B8").Value = "Modifing the lab_test table and add permanent value too it."
$ws.Range("A8").Value = "Redesign the database."
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "Md Mostafizur Rahman"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1

# Reselect cell F8 with row 7 scrolled to the top, matching the user's
# on-screen state after entering the new row.
$ws.Range("F8").Select()
$excel.ActiveWindow.ScrollRow = 7
